# "Script DataFrames Auxiliares Creado"
# Append the 2021 PAMI invoice rows (already present on sheet "Facturas PAMI 2021")
# to the bottom of sheet "Facturas PAMI 2020", and update the views/column widths
# to reflect where the user ended up working.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Facturas PAMI 2020"
$ws2 = $wb.Worksheets.Item(2)   # "Facturas PAMI 2021"

# ------------------------------------------------------------------
# 1) Copy the number formatting (styles) from the source rows on
#    "Facturas PAMI 2021" (A2:H9) onto the new destination rows
#    (A23:H30) on "Facturas PAMI 2020" before writing any values, so
#    the existing style records in the workbook get reused.
# ------------------------------------------------------------------
$ws2.Range("A2:A9").Copy() | Out-Null
$ws1.Range("A23").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

$ws2.Range("C2:H9").Copy() | Out-Null
$ws1.Range("C23").PasteSpecial(-4122) | Out-Null      # xlPasteFormats

# ------------------------------------------------------------------
# 2) Fill in the values / formulas for the new rows 23-30.
# ------------------------------------------------------------------
$rows = @(
    @{ A = 44200; B = 13222;  C = 78444711.159999996;  Cf = $null;                         E = 2526120.0499999998; F = 0;        G = -227540.65;             Gf = $null;                      H = $null;                  Hf = "=-9696292.58-6133.49" },
    @{ A = 44228; B = 13916;  C = 82646023.120000005;  Cf = $null;                         E = 999772.3;           F = 1715.7;   G = $null;                  Gf = "=-204072.27-4827.88";      H = $null;                  Hf = "=-7664309.47-595596.47" },
    @{ A = 44256; B = 14756;  C = 84682969.269999996;  Cf = $null;                         E = 918376.51;          F = 1205.72;  G = -192486.41;             Gf = $null;                      H = -6319174.9199999999;    Hf = $null },
    @{ A = 44292; B = 15523;  C = 85617373.560000002;  Cf = $null;                         E = 1505471.18;         F = 18200.66; G = -273730.84999999998;    Gf = $null;                      H = -11098717.789999999;    Hf = $null },
    @{ A = 44319; B = 16334;  C = 87512251.659999996;  Cf = $null;                         E = 1484598.43;         F = 5778.78;  G = -233094.29;             Gf = $null;                      H = -9631081.1500000004;    Hf = $null },
    @{ A = 44348; B = 88196;  C = $null;               Cf = "=77514763.83+21243509.85";    E = 2239762.44;         F = 1041.93;  G = -245940.37;             Gf = $null;                      H = -8831573.9100000001;    Hf = $null },
    @{ A = 44378; B = 156474; C = 106853801.95;        Cf = $null;                         E = 4347264.2300000004; F = 1890.9;   G = -280692.21000000002;    Gf = $null;                      H = -10129544.09;           Hf = $null },
    @{ A = 44411; B = 421851; C = $null;               Cf = "=83853447.67+21686755.93";    E = 3230463.38;         F = 9028.9;   G = -452225.06;             Gf = $null;                      H = -20183968.82;           Hf = $null }
)

$r = 23
foreach ($row in $rows) {
    $ws1.Cells.Item($r, 1).Value = $row.A
    $ws1.Cells.Item($r, 2).Value = $row.B

    if ($row.Cf) { $ws1.Cells.Item($r, 3).Formula = $row.Cf } else { $ws1.Cells.Item($r, 3).Value = $row.C }

    $ws1.Cells.Item($r, 4).Formula = "=C$r-SUM(E$r" + ":H$r)"

    $ws1.Cells.Item($r, 5).Value = $row.E
    $ws1.Cells.Item($r, 6).Value = $row.F

    if ($row.Gf) { $ws1.Cells.Item($r, 7).Formula = $row.Gf } else { $ws1.Cells.Item($r, 7).Value = $row.G }
    if ($row.Hf) { $ws1.Cells.Item($r, 8).Formula = $row.Hf } else { $ws1.Cells.Item($r, 8).Value = $row.H }

    $r = $r + 1
}

# Row 30 (E column) carries the "Verdana 8" cell style instead of the
# regular currency style - match that explicitly.
$ws1.Cells.Item(30, 5).Font.Name = "Verdana"
$ws1.Cells.Item(30, 5).Font.Size = 8

# ------------------------------------------------------------------
# 3) Column widths on "Facturas PAMI 2020": C and D used to share one
#    <col> entry (auto best-fit); now they are sized independently.
# ------------------------------------------------------------------
$ws1.Columns.Item(3).ColumnWidth = 20.167
$ws1.Columns.Item(4).ColumnWidth = 19.667

# ------------------------------------------------------------------
# 4) Update sheet selections: the user ended up with "Facturas PAMI
#    2020" selected near the bottom of the new data, while "Facturas
#    PAMI 2021" keeps a full-range selection but is no longer the
#    active tab.
# ------------------------------------------------------------------
$ws2.Range("A2:H9").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("A31").Select() | Out-Null
